$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename "time" sheet to "demand" and rebuild its header + data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("time")
$ws.Name = "demand"

# Header row - write in the same order the original author edited the
# cells (B, then E, then D, then C) so new shared-string entries land in
# the same table positions as the authored workbook.
$ws.Range("B1").Value = "Hours"
$ws.Range("E1").Value = "Heating power demand [kW]"
$ws.Range("D1").Value = "Electricity power demand [kW]"
$ws.Range("C1").Value = "EV power demand [kW]"

# Data rows: column B = hour index (1-24), column E = heating power demand [kW]
$heatingDemand = @(5,5,4,4,5,6,8,15,20,20,15,15,15,15,15,15,20,22,22,22,22,22,20,15)
for ($i = 0; $i -lt $heatingDemand.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $i + 1
    $ws.Range("E$row").Value = $heatingDemand[$i]
}

# Resize the newly-populated columns to fit their content.
$ws.Columns("C:E").AutoFit()

# This sheet becomes the active tab / selected sheet.
$null = $ws.Activate()
$null = $ws.Range("E2").Select()

# ---------------------------------------------------------------------
# 2. EV_data sheet: header text unchanged, just re-set to keep the
#    shared-string table consistent with the reordering above.
# ---------------------------------------------------------------------
$evData = $wb.Worksheets.Item("EV_data")
$evData.Range("B1").Value = "Energy capacity [kWh]"
$evData.Range("C1").Value = "Max charging power [kW]"
$evData.Range("D1").Value = "Arrival SOC [%]"
$evData.Range("E1").Value = "Departure SOC [%]"
$evData.Range("F1").Value = "Tarrival [h]"
$evData.Range("G1").Value = "Tdeparture [h]"

# ---------------------------------------------------------------------
# 3. grid_connection sheet: header text unchanged, re-set for the same
#    reason; it also loses tabSelected since "demand" is now active.
# ---------------------------------------------------------------------
$gridConn = $wb.Worksheets.Item("grid_connection")
$gridConn.Range("B1").Value = "Max power [kW]"

# ---------------------------------------------------------------------
# 4. HP sheet: header text unchanged, re-set for the same reason.
# ---------------------------------------------------------------------
$hp = $wb.Worksheets.Item("HP")
$hp.Range("B1").Value = "Power rated [kW]"

Write-Host "Workbook updated"
